# Refactor to make database structure template-based
#
# Adds "Contributor" / "Contact" columns (with their Instruction-sheet
# descriptions), renames several headers on "List of lines" to their
# display-cased form backed by the Instruction sheet, moves the "Note"
# column two slots to the right to make room, and makes "List of lines"
# the active tab/sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List of lines")
$ws2 = $wb.Worksheets.Item("Instruction")

# ---------------------------------------------------------------------
# 1) "List of lines": shift the trailing "Note" column from K to M and
#    insert the new "Contact"/"Citation" columns at J..L. Work from the
#    rightmost column backwards so we never clobber a value we still
#    need to copy.
# ---------------------------------------------------------------------

# Row 1 (headers) - move Note (K1) -> M1 first.
$ws1.Range("M1").Value = $ws1.Range("K1").Value

# Rows 2-3 ("This is an example..." note) - move K -> M.
$ws1.Range("M2").Value = $ws1.Range("K2").Value
$ws1.Range("M3").Value = $ws1.Range("K3").Value
$ws1.Range("K2").ClearContents()
$ws1.Range("K3").ClearContents()

# New headers for the inserted columns.
$ws1.Range("J1").Value = "Contributor"
$ws1.Range("K1").Value = "Contact"
$ws1.Range("L1").Value = "Citation"

# Re-cased headers (now backed by the Instruction sheet's own display text).
$ws1.Range("F1").Value = "Cassette style"
$ws1.Range("G1").Value = "Dimerization Domain"
$ws1.Range("H1").Value = "Status"
$ws1.Range("I1").Value = "Private"

# Column width: the 37.5-wide "Note" column formatting moves from K to M.
$ws1.Range("M1").ColumnWidth = 36.67
$ws1.Range("K1").ColumnWidth = 9.67

# A lightweight (no dropdown restriction) data validation now covers the
# re-cased G1:I1 header cells (just input-message placeholders).
$gi = $ws1.Range("G1:I1")
$gi.Validation.Add(0)
$gi.Validation.IgnoreBlank = $false

# ---------------------------------------------------------------------
# 2) "Instruction": add the matching Contributor/Contact columns (J, K)
#    with their header + description, and shift the old "citation"
#    header that used to live in J1 out to L1.
# ---------------------------------------------------------------------

$ws2.Range("L1").Value = $ws2.Range("J1").Value
$ws2.Range("J1").Value = "Contributor"
$ws2.Range("K1").Value = "Contact"

$ws2.Range("J2").Value = "Please leave blank if your lab generated this line. If you obtained this line from other researchers, please leave their name here."
$ws2.Range("K2").Value = "If left blank, the contact information from your account will be used. If you want people to request via a different mail, or if the line is from another lab and you know their preferred mail, please leave that here."
$ws2.Range("J2").WrapText = $true
$ws2.Range("K2").WrapText = $true

# Row 2 grows taller to fit the extra wrapped description text.
$ws2.Rows(2).RowHeight = 136

# Columns I..K (9-11) share the same 25-wide formatting.
$ws2.Range("I1:K1").ColumnWidth = 24.17

# ---------------------------------------------------------------------
# 3) View state: "List of lines" becomes the selected/active tab
#    (selection resting on J2); "Instruction" keeps a selection on K3
#    but is no longer the active tab.
# ---------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("K3").Select()

$ws1.Activate()
$ws1.Range("J2").Select()

Write-Output "done"
